$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark that sits in its own empty
#    paragraph right after the title blurb. Removing the (empty)
#    bookmark leaves a bare empty paragraph behind, matching the
#    target <w:p/>.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Locate the "October 2017" update-history paragraph and append,
#    after it, a blank paragraph followed by the new "November 2017"
#    paragraph describing the GetGamingDeviceModelInformation update.
# ------------------------------------------------------------------
$octPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "October 2017:*") {
        $octPara = $p
    }
}

$afterOct = $octPara.Range
$afterOct.Collapse(0)
$afterOct.InsertParagraphAfter()

# The paragraph that now sits right after "October 2017..." is the
# freshly inserted blank one; grab it and add another paragraph break
# after it so we end up with: October-2017 / <blank> / November-2017.
$blankPara = $octPara.Next()
$afterBlank = $blankPara.Range
$afterBlank.Collapse(0)
$afterBlank.InsertParagraphAfter()

$novPara = $blankPara.Next()
$novRange = $novPara.Range

$novRange.InsertAfter("November 2017: Updated for ")
$novRange.InsertAfter("GetGamingDeviceModelInformation")
$novRange.InsertAfter(" in the Fall Creators Update (16299)")
$novRange.InsertAfter(".")

# ------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark (collapsed/zero-length) at the
#    very end of the text we just typed, mirroring where Word leaves
#    it after the last edit.
# ------------------------------------------------------------------
$endOfNov = $d.Range($novRange.End, $novRange.End)
$d.Bookmarks.Add("_GoBack", $endOfNov)
